# Insert a new data row at row 211 of the "Papa" price sheet, shifting the
# existing rows 211..298 down to 212..299, then populate the new row with
# the latest weekly price observation (Feria Lagunitas de Puerto Montt,
# Pehuenche, "1a nueva(o)").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(211).Insert()

$ws.Range("A211").Value = 4
$ws.Range("B211").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C211").Value = "Los Lagos"
$ws.Range("D211").Value = 44553
$ws.Range("E211").Value = 10
$ws.Range("F211").Value = 100114001
$ws.Range("G211").Value = "Papa"
$ws.Range("H211").Value = "Pehuenche"
$ws.Range("I211").Value = "1a nueva(o)"
$ws.Range("J211").Value = 300
$ws.Range("K211").Value = 11000
$ws.Range("L211").Value = 12000
$ws.Range("M211").Value = 11500
$ws.Range("N211").Value = "`$/saco 25 kilos"
$ws.Range("O211").Value = "Región de La Araucanía"
$ws.Range("P211").Value = 460
$ws.Range("Q211").Value = 25
$ws.Range("R211").Value = "Hortaliza"
